$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rework the wording of the discussion-question paragraph (paragraph 2).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("say we ‘vectorize", $true, $false, $false, $false, $false, $true, 1, $false, "‘vectorize", 2) | Out-Null
$d.Content.Find.Execute("How do you think you can apply", $true, $false, $false, $false, $false, $true, 1, $false, "How can you apply", 2) | Out-Null
$d.Content.Find.Execute("one paragraph and you need to respond", $true, $false, $false, $false, $false, $true, 1, $false, "one paragraph, and you must respond", 2) | Out-Null
$d.Content.Find.Execute("student’s posts", $true, $false, $false, $false, $false, $true, 1, $false, "students’ posts", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Register the "Normal (Web)" paragraph style used by the pasted citation.
# ---------------------------------------------------------------------------
$webStyle = $d.Styles.Add("Normal (Web)", 1)
$webStyle.BaseStyle = "Normal"
$webStyle.Priority = 99
$webStyle.UnhideWhenUsed = $true

$webFmt = $webStyle.ParagraphFormat
$webFmt.SpaceBefore = 5
$webFmt.SpaceBeforeAuto = $true
$webFmt.SpaceAfter = 5
$webFmt.SpaceAfterAuto = $true
$webFmt.LineSpacingRule = 0

$webFont = $webStyle.Font
$webFont.Name = "Times New Roman"
$webFont.NameFarEast = "Times New Roman"
$webFont.NameBi = "Times New Roman"
$webFont.Size = 12
$webFont.SizeBi = 12

# ---------------------------------------------------------------------------
# 3. Append the new paragraphs after the discussion question:
#      - blank line
#      - the written answer
#      - "Sources:" label
#      - the citation (styled "Normal (Web)", hanging indent)
#      - trailing blank line
# ---------------------------------------------------------------------------
$q = $d.Paragraphs(2)
$q.Range.InsertParagraphAfter()

$blank1 = $d.Paragraphs(3)
$blank1.Range.InsertParagraphAfter()

$answerPara = $d.Paragraphs(4)
$answerPara.Range.Text = "A vectorized operation in R refers to an operation that applies to an entire vector as a single entity instead of having to act on each element individually. Most functions in R are vectorized, including arithmetic, comparison, and logical operators. One example of a vectorized function in use is the multiplication of a vector. For example, if you wanted to know what 1 through 10 times 2 was, you could create a vector containing 1 through 10. You could then multiply that vector by 2, which would output 2, 4, 6, 8, etc. This would be much more efficient than computing each element separately as 1 * 2, 2 * 2, etc. Overall, vectorization is an excellent time-saving skill in R and helps increase code efficiency and conciseness. "
$answerPara = $d.Paragraphs(4)
$answerPara.Range.InsertParagraphAfter()

$sourcesPara = $d.Paragraphs(5)
$sourcesPara.Range.Text = "Sources: "
$sourcesPara = $d.Paragraphs(5)
$sourcesPara.Range.InsertParagraphAfter()

$citationPara = $d.Paragraphs(6)
$citationPara.Range.Text = "Yale University. (n.d.). R for Novices: Vectorization. https://docs.ycrc.yale.edu/r-novice-gapminder/09-vectorization/#:~:text=Most%20of%20R’s%20functions%20are,read%2C%20and%20less%20error%20prone. "
$citationPara = $d.Paragraphs(6)
$citationPara.Range.Style = "Normal (Web)"
$citationPara = $d.Paragraphs(6)
$citationPara.Format.LeftIndent = 28.35
$citationPara.Format.FirstLineIndent = -28.35
$citationPara = $d.Paragraphs(6)
$citationPara.Range.InsertParagraphAfter()

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
